$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-x-IMU3-SA-A2")

# Replace the "180R" comment with "470R" for the R1, R2, R3 resistors row (C5).
# Leading apostrophe preserves the existing quote-prefix text formatting of the cell.
$ws.Range("C5").Value = "'470R"
